$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "E7";   Value = 15.1093 },
    @{ Cell = "D8";   Value = -8.746600000000006 },
    @{ Cell = "D10";  Value = -9.158599999999989 },
    @{ Cell = "D12";  Value = -7.239400000000001 },
    @{ Cell = "E14";  Value = 17.16429999999999 },
    @{ Cell = "E15";  Value = 15.7927 },
    @{ Cell = "D18";  Value = -8.566500000000001 },
    @{ Cell = "E18";  Value = 16.40120000000002 },
    @{ Cell = "E20";  Value = 15.6026 },
    @{ Cell = "D25";  Value = -8.145399999999999 },
    @{ Cell = "E29";  Value = 17.18430000000002 },
    @{ Cell = "E30";  Value = 15.65759999999999 },
    @{ Cell = "E31";  Value = 16.0193 },
    @{ Cell = "E35";  Value = 15.6401 },
    @{ Cell = "D37";  Value = -7.984099999999995 },
    @{ Cell = "E40";  Value = 17.04750000000002 },
    @{ Cell = "E44";  Value = 16.37209999999999 },
    @{ Cell = "E50";  Value = 16.3133 },
    @{ Cell = "E54";  Value = 16.47759999999999 },
    @{ Cell = "D55";  Value = -8.407100000000002 },
    @{ Cell = "D68";  Value = -7.776500000000003 },
    @{ Cell = "E68";  Value = 16.03779999999999 },
    @{ Cell = "E76";  Value = 16.26179999999999 },
    @{ Cell = "D77";  Value = -6.169199999999997 },
    @{ Cell = "D78";  Value = -7.503200000000006 },
    @{ Cell = "D79";  Value = -6.066800000000001 },
    @{ Cell = "D80";  Value = -7.492299999999996 },
    @{ Cell = "D81";  Value = -7.683199999999997 },
    @{ Cell = "D82";  Value = -8.443399999999995 },
    @{ Cell = "D84";  Value = -8.776100000000003 },
    @{ Cell = "E87";  Value = 16.37769999999999 },
    @{ Cell = "E88";  Value = 16.33409999999999 },
    @{ Cell = "E92";  Value = 18.42590000000003 },
    @{ Cell = "E96";  Value = 15.93059999999998 },
    @{ Cell = "E98";  Value = 15.28619999999999 },
    @{ Cell = "D101"; Value = -7.868600000000002 },
    @{ Cell = "E101"; Value = 16.663 },
    @{ Cell = "D102"; Value = -7.999300000000001 },
    @{ Cell = "E102"; Value = 16.5847 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
